# Remove the trailing "empty / page-break / copyright" block that used to
# follow the "LOM3008: Eletrometalurgia (Requisito fraco)" requirement line:
#   - an empty paragraph
#   - an empty page-break paragraph
#   - the "(c) 2020 . Contact: ..." paragraph
# leaving the final empty / page-break paragraph pair untouched.

$d = $word.ActiveDocument

# Locate the copyright paragraph ("(c) 2020 . Contact: luizeleno@usp.br...").
$copyrightIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
        break
    }
}

if ($copyrightIndex -eq -1) {
    throw "Could not find the copyright paragraph"
}

# The two empty paragraphs immediately preceding it (blank, then the
# page-break blank) are removed along with it.
$firstIndex = $copyrightIndex - 2

$startRange = $d.Paragraphs.Item($firstIndex).Range.Start
$endRange = $d.Paragraphs.Item($copyrightIndex).Range.End

$deleteRange = $d.Range($startRange, $endRange)
$deleteRange.Delete()
